$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- HKD section -----------------------------------------------------
# Insert a new deposit row at row 5 (pushes existing rows 5+ down by one).
$ws.Rows("5:5").Insert()

# New row 5 data.
$d1 = Get-Date -Year 2018 -Month 12 -Day 14 -Hour 0 -Minute 0 -Second 0
$d2 = Get-Date -Year 2019 -Month 12 -Day 16 -Hour 0 -Minute 0 -Second 0
$ws.Range("A5").Value = $d1
$ws.Range("B5").Value = $d2
$ws.Range("C5").Value = 20
$ws.Range("D5").Value = 2.5
$ws.Range("E5").Value = "* No uplift"
$ws.Range("G5").Formula = "=C5*D5*10000/100/12"

# Updated interest rates on the (now shifted) existing HKD rows.
$ws.Range("D2").Value = 2.35
$ws.Range("D3").Value = 2.35
$ws.Range("D6").Value = 2.55
$ws.Range("D7").Value = 2.55
$ws.Range("D8").Value = 2.35

# The "sum" row's C total historically missed the last data row; correct
# it now that an extra row has been inserted above it.
$ws.Range("C9").Formula = "=SUM(C2:C8)"

# --- Fund dividends ----------------------------------------------------
$ws.Range("G12").Value = 867.68

# --- USD section ---------------------------------------------------
$ws.Range("D18").Value = 2.8

# --- CNY section -----------------------------------------------------
# The uplift rate constant moved from 1.131 to 1.135 for all three rows.
$ws.Range("C23").Formula = "=E23*1.135"
$ws.Range("C24").Formula = "=E24*1.135"

# Replace the final CNY deposit row's data with the new entry.
$d3 = Get-Date -Year 2018 -Month 12 -Day 31 -Hour 0 -Minute 0 -Second 0
$d4 = Get-Date -Year 2019 -Month 12 -Day 31 -Hour 0 -Minute 0 -Second 0
$ws.Range("A25").Value = $d3
$ws.Range("B25").Value = $d4
$ws.Range("D25").Value = 3.6
$ws.Range("E25").Value = 27.5
$ws.Range("C25").Formula = "=E25*1.135"

# --- Restore the selection shown in the saved workbook --------------
$ws.Range("E9").Select()
